$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.688.00"
$ws.Range("D2").Style = 'Normal'
$ws.Range("E2").Value = "'  -0.02%  "
$ws.Range("E2").Style = 'Normal'
$ws.Range("D3").Value = "'1.647.06"
$ws.Range("D3").Style = 'Normal'
$ws.Range("E3").Value = "'  +0.63%  "
$ws.Range("E3").Style = 'Normal'
$ws.Range("E4").Value = "'  +0.31%  "
$ws.Range("E4").Style = 'Normal'
$ws.Range("D5").Value = "'216.20"
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = "'  +1.20%  "
$ws.Range("E5").Style = 'Normal'
$ws.Range("D6").Value = "'0.503"
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = "'  -0.43%  "
$ws.Range("E6").Style = 'Normal'
$ws.Range("E7").Value = "'  +0.30%  "
$ws.Range("E7").Style = 'Normal'
$ws.Range("D8").Value = "'0.253"
$ws.Range("D8").Style = 'Normal'
$ws.Range("E8").Value = "'  -0.49%  "
$ws.Range("E8").Style = 'Normal'
$ws.Range("E9").Value = "'  +0.54%  "
$ws.Range("E9").Style = 'Normal'
$ws.Range("D10").Value = "'19.36"
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = "'  +0.21%  "
$ws.Range("E10").Style = 'Normal'
$ws.Range("D11").Value = "'0.0844"
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = "'  +0.05%  "
$ws.Range("E11").Style = 'Normal'
$ws.Range("D12").Value = "'1.877.18"
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = "'  +0.64%  "
$ws.Range("E12").Style = 'Normal'
$ws.Range("E13").Value = "'  +3.24%  "
$ws.Range("E13").Style = 'Normal'
$ws.Range("D14").Value = "'1.643.52"
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = "'  +0.60%  "
$ws.Range("E14").Style = 'Normal'
$ws.Range("D15").Value = "'0.535"
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = "'  +1.45%  "
$ws.Range("E15").Style = 'Normal'
$ws.Range("D16").Value = "'66.28"
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = "'  +4.34%  "
$ws.Range("E16").Style = 'Normal'
$ws.Range("D17").Value = "'26.758.23"
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = "'  +0.25%  "
$ws.Range("E17").Style = 'Normal'
$ws.Range("E18").Value = "'  +1.38%  "
$ws.Range("E18").Style = 'Normal'
$ws.Range("D19").Value = "'219.67"
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = "'  -0.56%  "
$ws.Range("E19").Style = 'Normal'
$ws.Range("E20").Value = "'  +0.24%  "
$ws.Range("E20").Style = 'Normal'
$ws.Range("D21").Value = "'4.39"
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = "'  +1.72%  "
$ws.Range("E21").Style = 'Normal'
$ws.Range("D22").Value = "'6.33"
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = "'  +2.09%  "
$ws.Range("E22").Style = 'Normal'
$ws.Range("D23").Value = "'9.57"
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = "'  +1.15%  "
$ws.Range("E23").Style = 'Normal'
$ws.Range("E24").Value = "'  +9.97%  "
$ws.Range("E24").Style = 'Normal'
$ws.Range("D25").Value = "'147.17"
$ws.Range("D25").Style = 'Normal'
$ws.Range("E25").Value = "'  -0.73%  "
$ws.Range("E25").Style = 'Normal'
$ws.Range("E27").Value = "'  -0.71%  "
$ws.Range("E27").Style = 'Normal'
$ws.Range("D28").Value = "'7.11"
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = "'  +2.24%  "
$ws.Range("E28").Style = 'Normal'
$ws.Range("E29").Value = "'  +2.40%  "
$ws.Range("E29").Style = 'Normal'
$ws.Range("E30").Value = "'  +1.47%  "
$ws.Range("E30").Style = 'Normal'
$ws.Range("E31").Value = "'  +0.65%  "
$ws.Range("E31").Style = 'Normal'
$ws.Range("E32").Value = "'  +2.06%  "
$ws.Range("E32").Style = 'Normal'
$ws.Range("E33").Value = "'  +2.48%  "
$ws.Range("E33").Style = 'Normal'
$ws.Range("D34").Value = "'1.288.14"
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = "'  +5.81%  "
$ws.Range("E34").Style = 'Normal'
$ws.Range("E35").Value = "'  +1.72%  "
$ws.Range("E35").Style = 'Normal'
$ws.Range("D36").Value = "'0.0184"
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = "'  +6.45%  "
$ws.Range("E36").Style = 'Normal'
$ws.Range("E37").Value = "'  +0.42%  "
$ws.Range("E37").Style = 'Normal'
$ws.Range("E39").Value = "'  +2.00%  "
$ws.Range("E39").Style = 'Normal'
$ws.Range("D40").Value = "'1.00"
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = "'  +0.31%  "
$ws.Range("E40").Style = 'Normal'
$ws.Range("D42").Value = "'2.24"
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = "'  -1.92%  "
$ws.Range("E42").Style = 'Normal'
$ws.Range("E43").Value = "'  +0.00%  "
$ws.Range("E43").Style = 'Normal'
$ws.Range("D44").Value = "'1.789.30"
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = "'  +0.83%  "
$ws.Range("E44").Style = 'Normal'
$ws.Range("D45").Value = "'93.78"
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = "'  +0.59%  "
$ws.Range("E45").Style = 'Normal'
$ws.Range("D46").Value = "'60.22"
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = "'  +9.50%  "
$ws.Range("E46").Style = 'Normal'
$ws.Range("D47").Value = "'1.61"
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = "'  +3.30%  "
$ws.Range("E47").Style = 'Normal'
$ws.Range("D48").Value = "'0.0516"
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = "'  +0.65%  "
$ws.Range("E48").Style = 'Normal'
$ws.Range("D49").Value = "'7.82"
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = "'  +2.30%  "
$ws.Range("E49").Style = 'Normal'
$ws.Range("D50").Value = "'0.0978"
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = "'  +2.96%  "
$ws.Range("E50").Style = 'Normal'
$ws.Range("D51").Value = "'0.407"
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = "'  -0.79%  "
$ws.Range("E51").Style = 'Normal'
